$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the Celltype/Color table from columns F:G to columns A:B ---
$labels = @("Celltype","CD4 naive","CD4 TCM","CD4 TEM","CD4 CTL","CD4 Treg","CD4 proliferating","CD8 naive","CD8 TCM","CD8 TEM","CD8 CTL","CD8 proliferating","MAIT","NKT","DNT","GDT","NK CD16-","NK CD16+","NK proliferating","B naive","B intermediate","B memory","Plasma","Classical monocyte","Intermediate monocyte","Non-classical monocyte","CDC1","CDC2","PDC","ASDC","Platelet","HSPC","Lin-","Multiplet","Dead/debris")
$colors = @("Color","#7FC97F","#BEAED4","#FDC086","#FFFF99","#386CB0","#F0027F","#BF5B17","#E5D8BD","#1B9E77","#D95F02","#7570B3","#E7298A","#66A61E","#E6AB02","#A6761D","#666666","#A6CEE3","#1F78B4","#B2DF8A","#33A02C","#FB9A99","#E31A1C","#FDBF6F","#FF7F00","#CAB2D6","#6A3D9A","#FFFF99","#B15928","#FBB4AE","#B3CDE3","#CCEBC5","#DECBE4","#FED9A6","#FFFFCC")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value2 = $labels[$i]
    $ws.Cells.Item($r, 2).Value2 = $colors[$i]
}

# Clear out the old F:G columns now that the data lives in A:B
$ws.Range("F1:G35").ClearContents()

# --- Update the saved view/selection state ---
$ws.Range("F12").Select()
